# Apply updated dSF (column F) values for diekman_jake.xlsx
# Commit message: "repull data, push all data, mean calculation"
# This updates several rows in column F with re-pulled data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0
    8  = 5
    12 = -5
    13 = -4
    14 = -1
    15 = 0
    20 = -9
    23 = -1
    28 = -1
    35 = 0
    37 = 3
    41 = 0
    45 = 4
    50 = 0
    51 = -1
    54 = 2
    57 = -1
    64 = 2
    66 = 4
    68 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
